# Apply the changes described by the diff:
#  - Version bump: 1.0 -> 1.2.5
#  - Typo/punctuation fixes to several repeated precondition / description strings

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Version: 1.0 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# 2) Precondition text fix (appears in several test case blocks)
$oldPrecondition = "O usuario devidamente autenticado e na tela inicial do sistema"
$newPrecondition = "O usuário devidamente autenticado e na tela inicial do sistema."
$preconditionCells = @("B8", "B16", "B24", "B34", "B43", "B52", "B61", "B68", "B77", "B86")
foreach ($addr in $preconditionCells) {
    $ws.Range($addr).Value = $newPrecondition
}

# 3) TC1 step 2 expected result: fix "Permite não permite" -> "Não permite"
$ws.Range("D11").Value = "SYSTEM Identifica que a prestação de contas indicada pelo usuário não está em nenhum desses dois estados: a) NÃO REALIZADA e b) DEVOLVIDA; Não permite um novo envio ou alterações na prestação (exclusão de documentos)."

# 4) TC4 step 3 expected result: add trailing period
$ws.Range("D38").Value = "SYSTEM Apresenta a tela de Detalhar Diárias."

# 5) TC10 step 1: fix "histório" -> "histórico"
$ws.Range("B90").Value = "Chefe Verifica o histórico da tramitação da prestação de contas."
